# "new graphs for sups ratings"
# 1) Reorder sheet tabs: Ratings, Work Again, Additional Comments, Past Supervisors, Absent
# 2) Reorder the per-supervisor "Why work again" comments on the "Work Again" sheet
# 3) Reorder the per-supervisor notes on the "Additional Comments" sheet
# 4) Append two new rows (Lawrence, Godfrey) to "Past Supervisors"
# 5) Replace the "Absent" sheet's data with a smaller, consolidated table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Move "Absent" to the end of the tab strip, then re-select
#    "Additional Comments" so it stays the active tab.
# ---------------------------------------------------------------------------
$absent = $wb.Worksheets.Item("Absent")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$absent.Move($null, $lastSheet)

# ---------------------------------------------------------------------------
# 2) "Work Again" sheet -- reorder comments within each supervisor block.
# ---------------------------------------------------------------------------
$wa = $wb.Worksheets.Item("Work Again")

$wa.Range("B2").Value = "It portrayed provided  a good working enviroment  for the whole team."
$wa.Range("B3").Value = "She is patient, honest and understanding"
$wa.Range("B4").Value = "It would be easy for me to work with her because I know what she likes and dislikes."
$wa.Range("B5").Value = "She is approachable and a good mentor."
$wa.Range("B6").Value = "She is so supportive"
$wa.Range("B7").Value = "She is honest and passes clear information and she is approachable"

$wa.Range("B9").Value = "He easily understands situations when informed."
$wa.Range("B10").Value = "He is sociable and easy to approach ."
$wa.Range("B11").Value = "Because i keep learning"
$wa.Range("B12").Value = "He treated me as an individual,never at one one time did he assume things about me."

$wa.Range("B13").Value = "Because he was fair, clear in his communication and a good teacher"
$wa.Range("B14").Value = "Because hes freindly."
$wa.Range("B15").Value = "He is good"

$wa.Range("B16").Value = "You are not worried about failing to accomplish a challenge as she always finds a way out unless otherwise."
$wa.Range("B18").Value = "She has been really good to me and to the whole team i think"
$wa.Range("B19").Value = "She is hardworking and exemplary"
$wa.Range("B20").Value = "She is  polite and  honest"
$wa.Range("B22").Value = "She is harding and team player"

$wa.Range("B23").Value = "Because she has room for improvement in the weak areas."
$wa.Range("B24").Value = "She is fun to work with"
$wa.Range("B26").Value = "I wish to cope up her expression nd the way she does things."

$wa.Range("B31").Value = "Too bossy"
$wa.Range("B32").Value = "I don't want someone  who  intervenes in my privacy when it's not part  of field work."

# ---------------------------------------------------------------------------
# 3) "Additional Comments" sheet -- reorder each supervisor's remarks.
# ---------------------------------------------------------------------------
$ac = $wb.Worksheets.Item("Additional Comments")

$ac.Range("B2").Value = "She is a wonderful leader that i would wish to work with next time"
$ac.Range("B3").Value = "I think it's all said!"
$ac.Range("B4").Value = "She should improve on her communication about work plan"
$ac.Range("B5").Value = "She is good leader and good to work with"
$ac.Range("B6").Value = "She's very very very very very very good"

$ac.Range("B7").Value = "He is short tampered. However he is a direct and honest guy."
$ac.Range("B8").Value = "He should be considered for any future opportunities."
$ac.Range("B9").Value = "Good"
$ac.Range("B10").Value = "I think he should work on being fast about communicating."

$ac.Range("B11").Value = "She's good and still good to be a supervisor incase of another project."
$ac.Range("B12").Value = "She is the best"
$ac.Range("B13").Value = "She is good to go.....It would be my pleasure if she remains a supervisor"
$ac.Range("B14").Value = "She was a good supervisor"
$ac.Range("B15").Value = "She is good for me"
$ac.Range("B16").Value = "She is friendly and has managed to pool contacts of many OC traffic officers"
$ac.Range("B17").Value = "She is good and she loves what she does"

$ac.Range("B18").Value = "Just fair."
$ac.Range("B19").Value = "Figures stress Julie."
$ac.Range("B20").Value = "  a beleiver, just like me. I like her."

$ac.Range("B21").Value = "She is good"
$ac.Range("B22").Value = "Has a lot to improve"

# ---------------------------------------------------------------------------
# 4) "Past Supervisors" sheet -- append Lawrence / Godfrey rows.
# ---------------------------------------------------------------------------
$ps = $wb.Worksheets.Item("Past Supervisors")

$ps.Range("A9").Value = "Lawrence"
$ps.Range("B9").Value = 2
$ps.Range("A10").Value = "Godfrey"
$ps.Range("B10").Value = 1

# ---------------------------------------------------------------------------
# 5) "Absent" sheet -- replace with the smaller, consolidated table.
# ---------------------------------------------------------------------------
$ab = $wb.Worksheets.Item("Absent")
$ab.Range("A1:C8").ClearContents()

$ab.Range("A1").Value = "Are you aware of any incidents in which you had difficulty reaching out to your supervisor?"
$ab.Range("A2").Value = "How many times was your supervisor difficult to reach?"

$ab.Range("B4").Value = "Supervisor"
$ab.Range("C4").Value = "How Many Times Absent"
$ab.Range("B5").Value = "Isaac"
$ab.Range("C5").Value = 2
$ab.Range("B6").Value = "Julie"
$ab.Range("C6").Value = 2

# ---------------------------------------------------------------------------
# Finally, make sure "Additional Comments" is the selected/active tab.
# ---------------------------------------------------------------------------
$ac.Activate()
